$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume data per upstream scrape refresh.
# For cells whose new value is a plain decimal number (e.g. "603.34"),
# force text format first so Excel does not silently coerce it to a
# binary double (losing trailing zeros / exact decimal text) or switch
# tiny values to scientific notation. All these columns are plain text
# columns in the source sheet (t="inlineStr"), never real numbers.

$ws.Range('D2').Value = '68.578.05'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '3.903.96'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.34'
$ws.Range('E5').Value = '  +0.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.69'
$ws.Range('E6').Value = '  +2.13%  '
$ws.Range('D7').Value = '3.904.80'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.531'
$ws.Range('E9').Value = '  +0.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.47'
$ws.Range('E11').Value = '  +0.64%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.460'
$ws.Range('E12').Value = '  +0.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000255'
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.19'
$ws.Range('E14').Value = '  -0.39%  '
$ws.Range('D15').Value = '4.557.76'
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('D16').Value = '3.895.48'
$ws.Range('E16').Value = '  +0.45%  '
$ws.Range('D17').Value = '68.476.32'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.15'
$ws.Range('E18').Value = '  +5.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.43'
$ws.Range('E19').Value = '  -0.78%  '
$ws.Range('E20').Value = '  +0.30%  '
$ws.Range('E21').Value = '  -1.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '473.71'
$ws.Range('E22').Value = '  -2.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.741'
$ws.Range('E23').Value = '  +2.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000168'
$ws.Range('E24').Value = '  +0.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.83'
$ws.Range('E25').Value = '  -0.95%  '
$ws.Range('E26').Value = '  +1.42%  '
$ws.Range('E27').Value = '  +1.40%  '
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('E29').Value = '  -1.06%  '
$ws.Range('E30').Value = '  +1.43%  '
$ws.Range('D31').Value = '4.051.89'
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.88'
$ws.Range('E32').Value = '  +1.86%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '31.61'
$ws.Range('E33').Value = '  -0.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.32'
$ws.Range('E34').Value = '  -2.45%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.45'
$ws.Range('E35').Value = '  +1.47%  '
$ws.Range('D36').Value = '3.874.45'
$ws.Range('E36').Value = '  +0.35%  '
$ws.Range('E37').Value = '  -1.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.73'
$ws.Range('E38').Value = '  +16.59%  '
$ws.Range('E39').Value = '  -0.14%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.142'
$ws.Range('E40').Value = '  +2.15%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('E44').Value = '  +15.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '429.76'
$ws.Range('E45').Value = '  +0.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.01'
$ws.Range('E46').Value = '  +1.16%  '
$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.63'
$ws.Range('E48').Value = '  +1.32%  '
$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '47.28'
$ws.Range('E49').Value = '  -1.87%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '27.20'
$ws.Range('E50').Value = '  +6.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '144.17'
$ws.Range('E51').Value = '  +1.08%  '
